# Adds two new paragraphs of body text under the "Problem Statement"
# heading, in the empty paragraph that currently follows it.
#
# Paragraph 1 (fills the existing empty paragraph):
#   "The most significant sense in the human body is vision. ..."
#
# Paragraph 2 (brand-new paragraph inserted after it):
#   "For instance, ... application called "" + italic("Object Identifier") +
#   "" that detects objects ... mimic the human eye."

$d = $word.ActiveDocument
$cr = [string][char]13

$para1Text = 'The most significant sense in the human body is vision. It enables one to assess and comprehend one''s environment. At least 280 million individuals are visually challenged or visually impaired, according to data obtained from the WHO. One''s daily activities may be hampered by vision issues. Examples of these difficulties include reading text, crossing streets, and recognizing items in daily life. Although they can learn other coping mechanisms, they encounter certain navigational challenges and social awkwardness.'

$para2TextA = 'For instance, people have trouble locating a certain room in a strange setting. It can be challenging for persons who are visually impaired to tell during a discussion whether someone is speaking to them directly or to someone else. To assist those who are blind, a straightforward Android application called "'
$para2TextB = 'Object Identifier'
$para2TextC = '" that detects objects is the basis of the suggested solution. With the aid of a smartphone camera and object detection, this application aims to mimic the human eye.'

# The empty paragraph right after the "Problem Statement" heading.
$targetPara = $d.Paragraphs(4)
Write-Output ("target paragraph text before = [" + $targetPara.Range.Text + "]")

# Split it into three empty paragraphs by inserting two paragraph marks
# after its (empty) range. The original paragraph stays in slot 4; two
# brand-new paragraphs (with no rsid/paraId baggage, just like freshly
# authored content) land in slots 5 and 6.
$targetPara.Range.InsertAfter($cr + $cr)

# Paragraph 5: new -> becomes the first body paragraph.
$p5 = $d.Paragraphs(5)
$p5.Range.InsertBefore($para1Text)

# Paragraph 6: new -> becomes the second body paragraph, with an italic
# run in the middle.
$p6 = $d.Paragraphs(6)
$p6Start = $p6.Range.Start
$p6.Range.InsertBefore($para2TextA + $para2TextB + $para2TextC)

$italicStart = $p6Start + $para2TextA.Length
$italicEnd = $italicStart + $para2TextB.Length
$italicRange = $d.Range($italicStart, $italicEnd)
Write-Output ("italic run text = [" + $italicRange.Text + "]")
$italicRange.Italic = $true

# Remove the now-redundant original empty paragraph (slot 4), which
# merges the two new paragraphs up into slots 4 and 5.
$d.Paragraphs(4).Range.Delete()

Write-Output ("paragraph count = " + $d.Paragraphs.Count)
foreach ($pp in $d.Paragraphs) {
    Write-Output ("[" + $pp.Range.Text + "]")
}
